$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verleiherabgaben")
$ws.Columns("B:B").Insert()
$ws.Range("D9").NumberFormat = '"CHF"\ #,##0.00'
